$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "245.21")
# are not auto-converted to numbers by Excel, matching the original inline-string text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.892.17'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.732.19'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '245.21'
$ws.Range("E5").Value = '  +2.96%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").Value = '0.5001'
$ws.Range("E7").Value = '  -2.93%  '
$ws.Range("D8").Value = '0.2708'
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("D9").Value = '0.06149'
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").Value = '1.736.68'
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").Value = '0.07238'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '0.6510'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").Value = '4.763'
$ws.Range("E14").Value = '  +3.61%  '
$ws.Range("D15").Value = '76.86'
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '25.892.65'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").Value = '0.000006795'
$ws.Range("D21").Value = '4.616'
$ws.Range("E21").Value = '  +8.06%  '
$ws.Range("D22").Value = '1.960.43'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '8.747'
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("D24").Value = '5.432'
$ws.Range("E24").Value = '  +3.46%  '
$ws.Range("D25").Value = '133.49'
$ws.Range("E25").Value = '  -3.84%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '1.775'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '1.403'
$ws.Range("E28").Value = '  -7.18%  '
$ws.Range("D29").Value = '105.36'
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").Value = '3.972'
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").Value = '0.08089'
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("D33").Value = '0.04729'
$ws.Range("E33").Value = '  +3.06%  '
$ws.Range("D34").Value = '2.658'
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = '0.9925'
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").Value = '0.6084'
$ws.Range("E36").Value = '  -1.98%  '
$ws.Range("D37").Value = '2.731'
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("D38").Value = '0.01600'
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("D39").Value = '0.8603'
$ws.Range("E39").Value = '  +16.98%  '
$ws.Range("D40").Value = '1.933'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = '100.01'
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").Value = '0.3899'
$ws.Range("E43").Value = '  +1.11%  '
$ws.Range("D44").Value = '5.001'
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("E45").Value = '  +4.17%  '
$ws.Range("D46").Value = '6.303'
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").Value = '55.50'
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '30.64'
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.620'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3468'
$ws.Range("E51").Value = '  +1.30%  '

# Restore column D to the default "Normal" style so no stray number formatting remains.
$ws.Range("D2:D51").Style = "Normal"

